# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) was previously populated from a different
# source field (Strike#) and is being regenerated here using the
# correct per-game strikeout ("K") values. We write the recalculated
# values directly into column G for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => new K value (s_vals), keyed by worksheet row number
$s_vals = [ordered]@{
    2  = 1
    3  = 2
    4  = 2
    5  = 1
    6  = 2
    7  = 2
    8  = 0
    9  = 1
    10 = 2
    11 = 2
    12 = 2
    13 = 0
    14 = 2
    15 = 0
    16 = 1
    17 = 0
    18 = 2
    19 = 0
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 2
    30 = 1
    31 = 0
    32 = 1
    33 = 2
    34 = 1
    35 = 1
    36 = 3
    37 = 2
    38 = 1
    39 = 0
    40 = 1
    41 = 1
    42 = 0
    43 = 1
    44 = 0
    45 = 2
    46 = 0
    47 = 1
    48 = 1
    49 = 2
    50 = 0
    51 = 3
    52 = 1
    53 = 1
    54 = 3
    55 = 1
}

# Column G holds "K" (strikeouts). Write the regenerated value for each row.
foreach ($row in $s_vals.Keys) {
    $ws.Cells.Item($row, 7).Value = $s_vals[$row]
}
